# Add a new worksheet ("Sheet2") after the last existing sheet (Sheet1),
# then write the sample text into cell B2 of the new sheet.
$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Sheet2"

$newSheet.Range("B2").Value = "Lorem ipsum dolor sit amet, consectetur adipiscing elit, sed do eiusmod tempor"
